# Localization status workbook: the "Status" value for this file moved from
# "Ready for handoff" to "In Translation" (zh-cn / de-de). The same string is
# shown on the Overview sheet (columns E/F) and on each per-locale sheet
# (column C). Updating the cell text lets Excel's shared-string table drop the
# now-unused "Ready for handoff" entry and add "In Translation" automatically.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "In Translation"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "In Translation"

# The Status column narrowed to fit the shorter text ("In Translation" is
# shorter than "Ready for handoff") on all three sheets that show it.
# (Column letters: Overview E/F, zh-cn/de-de column C -> numeric indices 5/6/3.)
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5
$zhcn.Columns.Item(3).ColumnWidth = 12.5
$dede.Columns.Item(3).ColumnWidth = 12.5
